# Auto update stock data
# Applies:
#  1) Refresh date/EBITDA headline values on the existing 6 sheets (2025/10/29 -> 2025/10/30, etc.)
#  2) Fill in Ryerson Holding's previously-blank Altman Z-Score column (2.75)
#  3) Append 7 new per-company sheets with their historical metrics

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a value to be stored as literal text (mirrors typing a value
# preceded by an apostrophe in Excel) so numeric-looking strings like
# "12.05" or date-looking strings like "2025/10/30" are NOT auto-converted
# to numbers / date serials. An apostrophe with nothing after it produces an
# explicitly-typed, empty text cell (matches the workbook's existing blank
# Altman Z-Score / Piotroski F-Score placeholder cells).
# ---------------------------------------------------------------------------
function Set-Text($range, [string]$text) {
    $range.Value = "'" + $text
}

# ===========================================================================
# 1) Update row 2 ("most recent") figures on the six existing sheets
# ===========================================================================

# Alcoa (AA)
$ws = $wb.Worksheets.Item(1)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "4.73"

# Rio Tinto (RIO)
$ws = $wb.Worksheets.Item(2)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "7.81"

# Norsk Hydro (NHY)
$ws = $wb.Worksheets.Item(3)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "2.71"

# Reliance Steel & Aluminum (RS)
$ws = $wb.Worksheets.Item(4)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "12.44"

# Kaiser Aluminum (KALU)
$ws = $wb.Worksheets.Item(5)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "9.76"

# Ryerson Holding (RYI)
$ws = $wb.Worksheets.Item(6)
Set-Text $ws.Cells.Item(2,1) "2025/10/30"
Set-Text $ws.Cells.Item(2,2) "25.73"
# Altman Z-Score column (G) was blank for every historical row; now populated
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r,7).Value = 2.75
}

# ===========================================================================
# 2) Append the seven new company sheets (sheetId 7-13), in order, right
#    after "Ryerson Holding"
# ===========================================================================

$headers8 = @("Date_1","EBITDA","Debt / Equity Ratio","Inventory Turnover","Current Ratio","Ticker","Altman Z-Score","Piotroski F-Score")
$headers5 = @("Date_1","Debt / Equity Ratio","Ticker","Altman Z-Score","Piotroski F-Score")

function Add-SheetAtEnd([string]$name) {
    $wbLocal = $excel.ActiveWorkbook
    $lastSheet = $wbLocal.Worksheets.Item($wbLocal.Worksheets.Count)
    $newWs = $wbLocal.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newWs.Name = $name
    return $newWs
}

function Write-HeaderRow($ws, [string[]]$headers) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i+1).Value = $headers[$i]
    }
}

# --- Ultra Clean Holdings (ULTR) -------------------------------------------------
$ws = Add-SheetAtEnd "Ultra Clean Holdings"
Write-HeaderRow $ws $headers8
$rows = @(
    @("2025/10/30","12.05","0.83","4.52","3.21","ULTR"),
    @("2024/12/31","12.24","0.71","4.61","2.89","ULTR"),
    @("2023/12/31","18.45","0.71","3.56","2.88","ULTR"),
    @("2022/12/31","6.44","0.65","4.64","2.82","ULTR"),
    @("2021/12/31","11.12","0.71","5.97","2.42","ULTR"),
    @("2020/12/31","8.71","0.57","6.27","2.71","ULTR"),
    @("2015/12/31","Upgrade","Upgrade","Upgrade","Upgrade","ULTR")
)
$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        Set-Text $ws.Cells.Item($r, $i+1) $row[$i]
    }
    Set-Text $ws.Cells.Item($r,7) ""
    Set-Text $ws.Cells.Item($r,8) ""
    $r++
}

# --- Foxconn (header only, no historical rows yet) --------------------------
$ws = Add-SheetAtEnd "Foxconn"
Write-HeaderRow $ws $headers8

# --- Ferrotec Holdings (header only, no historical rows yet) ----------------
$ws = Add-SheetAtEnd "Ferrotec Holdings"
Write-HeaderRow $ws $headers8

# --- Benchmark Electronics (BHE) --------------------------------------------
$ws = Add-SheetAtEnd "Benchmark Electronics"
Write-HeaderRow $ws $headers8
$rows = @(
    @("2025/10/30","11.24","0.30","4.11","2.36","BHE"),
    @("2024/12/31","11.01","0.34","3.86","2.32","BHE"),
    @("2023/12/31","8.13","0.44","3.64","2.29","BHE"),
    @("2022/12/31","7.88","0.41","4.21","2.21","BHE"),
    @("2021/12/31","8.74","0.25","4.82","2.06","BHE"),
    @("2020/12/31","10.03","0.24","5.85","2.50","BHE"),
    @("2015/12/31","Upgrade","Upgrade","Upgrade","Upgrade","BHE")
)
$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        Set-Text $ws.Cells.Item($r, $i+1) $row[$i]
    }
    $ws.Cells.Item($r,7).Value = 3.3
    $ws.Cells.Item($r,8).Value = 8
    $r++
}

# --- Celestica (CLS) ---------------------------------------------------------
$ws = Add-SheetAtEnd "Celestica"
Write-HeaderRow $ws $headers5
$rows = @(
    @("2025/10/30","0.18","CLS"),
    @("2024/12/31","0.12","CLS"),
    @("2023/12/31","0.23","CLS"),
    @("2022/12/31","0.10","CLS"),
    @("2021/12/31","0.09","CLS"),
    @("2020/12/31","0.18","CLS"),
    @("2019/12/31","Upgrade","CLS")
)
$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        Set-Text $ws.Cells.Item($r, $i+1) $row[$i]
    }
    Set-Text $ws.Cells.Item($r,4) ""
    $ws.Cells.Item($r,5).Value = 5
    $r++
}

# --- Flex Ltd (FLEX) ---------------------------------------------------------
$ws = Add-SheetAtEnd "Flex Ltd"
Write-HeaderRow $ws $headers8
$rows = @(
    @("2025/10/30","13.89","0.85","4.46","1.33","FLEX"),
    @("2025/12/31","8.14","0.86","4.18","1.30","FLEX"),
    @("2024/12/31","8.83","0.73","3.59","1.52","FLEX"),
    @("2023/12/31","8.04","0.76","3.80","1.48","FLEX"),
    @("2022/12/31","7.94","1.16","4.36","1.34","FLEX"),
    @("2021/12/31","7.56","1.30","5.82","1.45","FLEX"),
    @("2016/12/31","Upgrade","Upgrade","Upgrade","Upgrade","FLEX")
)
$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        Set-Text $ws.Cells.Item($r, $i+1) $row[$i]
    }
    Set-Text $ws.Cells.Item($r,7) ""
    $ws.Cells.Item($r,8).Value = 7
    $r++
}

# --- MKS Instruments (MKS) ---------------------------------------------------
$ws = Add-SheetAtEnd "MKS Instruments"
Write-HeaderRow $ws $headers8
$rows = @(
    @("2025/10/30","15.22","1.84","2.12","2.99","MKS"),
    @("2024/12/31","12.96","2.06","1.99","3.19","MKS"),
    @("2023/12/31","14.25","2.03","2.01","3.18","MKS"),
    @("2022/12/31","11.25","1.15","2.57","2.94","MKS"),
    @("2021/12/31","11.63","0.36","2.91","4.67","MKS"),
    @("2020/12/31","15.46","0.44","2.66","4.83","MKS"),
    @("2015/12/31","Upgrade","Upgrade","Upgrade","Upgrade","MKS")
)
$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        Set-Text $ws.Cells.Item($r, $i+1) $row[$i]
    }
    Set-Text $ws.Cells.Item($r,7) ""
    Set-Text $ws.Cells.Item($r,8) ""
    $r++
}

# Restore the selection/active sheet to the first sheet like the original file
$wb.Worksheets.Item(1).Select()
